$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("B1").Value = "FALSE_count"
$ws.Range("C1").Value = "FALSE_percent"
$ws.Range("D1").Value = "TRUE_count"
$ws.Range("E1").Value = "TRUE_percent"

# Row 2 values
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 66.66666666666666
$ws.Range("D2").Value = 50
$ws.Range("E2").Value = 67.56756756756756

# Row 3 values
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 22.22222222222222
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 9.45945945945946

# Row 4 values
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 11.11111111111111
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 22.97297297297298
